$wb = $excel.ActiveWorkbook

# Rename sheets (drop the "_Presence" / fix typo suffix)
$wb.Worksheets.Item(1).Name = "m9-EX_glc__D_e10.0-EX_glu__L_e8.33"
$wb.Worksheets.Item(2).Name = "m9-EX_glc__D_e10.0-EX_gly_e3.33"
$wb.Worksheets.Item(3).Name = "m9-EX_glc__D_e10.0-EX_thr__L_e6.66"

# ---- Sheet 1 data ----
$ws = $wb.Worksheets.Item(1)
$data = @(
    @(2, 17, "Cofactor and Prosthetic Group Biosynthesis", [double]"4.177335996661073e-35", 0.6666666666666666, 219),
    @(3, 25, "Membrane Lipid Metabolism", [double]"2.644452797142522e-14", 0.8536585365853658, 41),
    @(4, 0, "Purine and Pyrimidine Biosynthesis", [double]"1.005000077795777e-12", 0.96, 25),
    @(5, 30, "Valine, Leucine, and Isoleucine Metabolism", [double]"1.537091919039133e-09", 1, 16),
    @(6, 37, "Glycolysis/Gluconeogenesis", [double]"2.041469048564413e-08", 0.8636363636363636, 22),
    @(7, 6, "Tyrosine, Tryptophan, and Phenylalanine Metabolism", [double]"8.719332125373817e-07", 0.7826086956521739, 23),
    @(8, 2, "Histidine Metabolism", [double]"3.19615229056376e-06", 1, 10),
    @(9, 35, "Citric Acid Cycle", [double]"1.258660091336561e-05", 0.8571428571428571, 14),
    @(10, 3, "Pentose Phosphate Pathway", 0.000117150557472181, 0.8333333333333334, 12),
    @(11, 4, "Threonine and Lysine Metabolism", 0.0003168068471862367, 0.6842105263157895, 19),
    @(12, 27, "Murein Biosynthesis", 0.0003748206339453732, 0.7333333333333333, 15),
    @(13, 38, "Arginine and Proline Metabolism", 0.001006195905647336, 0.525, 40),
    @(14, 10, "Lipopolysaccharide Biosynthesis / Recycling", 0.001701022587909553, 0.4507042253521127, 71),
    @(15, 29, "Nucleotide Salvage Pathway", 0.002201852499269651, 0.3956834532374101, 139),
    @(16, 1, "Alanine and Aspartate Metabolism", 0.002945603979694493, 0.7777777777777778, 9),
    @(17, 36, "Methionine Metabolism", 0.01007911282027065, 0.6, 15),
    @(18, 7, "Cysteine Metabolism", 0.01253412044890904, 0.6153846153846154, 13),
    @(19, 23, "Pyruvate Metabolism", 0.03592929644152482, 0.6, 10),
    @(20, 13, "Intracellular demand", 0.05784426907355551, 0.6666666666666666, 6),
    @(21, 24, "Biomass and maintenance functions", 0.0715342018489301, 0.75, 4),
    @(22, 11, "Anaplerotic Reactions", 0.1647293636403392, 0.5, 8),
    @(23, 18, "Oxidative Phosphorylation", 0.1935539112072101, 0.3461538461538461, 52),
    @(24, 9, "Glutamate Metabolism", 0.2253670070064976, 0.5, 6),
    @(25, 12, "Folate Metabolism", 0.233008917336096, 0.4444444444444444, 9),
    @(26, 8, "Glycine and Serine Metabolism", 0.3616317263225998, 0.3571428571428572, 14),
    @(27, 20, "Cell Envelope Biosynthesis", 0.453807910678834, 0.291044776119403, 134),
    @(28, 32, "Inorganic Ion Transport and Metabolism", 0.7539340399193886, 0.2589285714285715, 112),
    @(29, 14, "Unassigned", 0.9895606643899191, 0.1153846153846154, 26),
    @(30, 5, "Murein Recycling", 0.9996290407812678, 0.07894736842105263, 38),
    @(31, 15, "Transport, Inner Membrane", 0.9999972302661364, 0.1837349397590362, 332),
    @(32, 31, "Alternate Carbon Metabolism", 0.9999997984882145, 0.1384615384615385, 195),
    @(33, 22, "Transport, Outer Membrane", 0.9999998124269528, 0.02173913043478261, 46),
    @(34, 21, "Glycerophospholipid Metabolism", 0.999999999998509, 0.1097560975609756, 246),
    @(35, 26, "Glyoxylate Metabolism", 1, 0, 4),
    @(36, 28, "Transport, Outer Membrane Porin", 1, 0.08888888888888889, 270),
    @(37, 16, "Extracellular exchange", 1, 0.07716049382716049, 324),
    @(38, 33, "Methylglyoxal Metabolism", 1, 0, 9),
    @(39, 34, "Nitrogen Metabolism", 1, 0, 13),
    @(40, 19, "tRNA Charging", 1, 0, 22)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(39, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(40, 1).PasteSpecial(-4122) | Out-Null

# ---- Sheet 2 data ----
$ws = $wb.Worksheets.Item(2)
$data = @(
    @(2, 17, "Cofactor and Prosthetic Group Biosynthesis", [double]"3.902060658646019e-36", 0.6712328767123288, 219),
    @(3, 25, "Membrane Lipid Metabolism", [double]"2.303927306096845e-14", 0.8536585365853658, 41),
    @(4, 0, "Purine and Pyrimidine Biosynthesis", [double]"2.916157336609939e-11", 0.92, 25),
    @(5, 30, "Valine, Leucine, and Isoleucine Metabolism", [double]"1.438354038623101e-09", 1, 16),
    @(6, 6, "Tyrosine, Tryptophan, and Phenylalanine Metabolism", [double]"6.024236325619492e-09", 0.8695652173913043, 23),
    @(7, 37, "Glycolysis/Gluconeogenesis", [double]"2.400889281817839e-07", 0.8181818181818182, 22),
    @(8, 2, "Histidine Metabolism", [double]"3.066769285430075e-06", 1, 10),
    @(9, 3, "Pentose Phosphate Pathway", 0.0001127271797358665, 0.8333333333333334, 12),
    @(10, 38, "Arginine and Proline Metabolism", 0.0003028695451739138, 0.55, 40),
    @(11, 27, "Murein Biosynthesis", 0.0003601995126364794, 0.7333333333333333, 15),
    @(12, 35, "Citric Acid Cycle", 0.0009556003024706312, 0.7142857142857143, 14),
    @(13, 4, "Threonine and Lysine Metabolism", 0.001521330823492825, 0.631578947368421, 19),
    @(14, 10, "Lipopolysaccharide Biosynthesis / Recycling", 0.00157841230174295, 0.4507042253521127, 71),
    @(15, 36, "Methionine Metabolism", 0.002146401882138933, 0.6666666666666666, 15),
    @(16, 18, "Oxidative Phosphorylation", 0.009237581514104472, 0.4423076923076923, 52),
    @(17, 7, "Cysteine Metabolism", 0.01220910762765804, 0.6153846153846154, 13),
    @(18, 23, "Pyruvate Metabolism", 0.0352316527863909, 0.6, 10),
    @(19, 9, "Glutamate Metabolism", 0.05703465230232828, 0.6666666666666666, 6),
    @(20, 13, "Intracellular demand", 0.05703465230232828, 0.6666666666666666, 6),
    @(21, 8, "Glycine and Serine Metabolism", 0.06935886365141529, 0.5, 14),
    @(22, 24, "Biomass and maintenance functions", 0.07073527754603325, 0.75, 4),
    @(23, 1, "Alanine and Aspartate Metabolism", 0.07814319336731958, 0.5555555555555556, 9),
    @(24, 12, "Folate Metabolism", 0.07814319336731958, 0.5555555555555556, 9),
    @(25, 11, "Anaplerotic Reactions", 0.1627825078924616, 0.5, 8),
    @(26, 26, "Glyoxylate Metabolism", 0.3170594046903736, 0.5, 4),
    @(27, 29, "Nucleotide Salvage Pathway", 0.3264550675093258, 0.302158273381295, 139),
    @(28, 32, "Inorganic Ion Transport and Metabolism", 0.4189024218698006, 0.2946428571428572, 112),
    @(29, 20, "Cell Envelope Biosynthesis", 0.4416369085763585, 0.291044776119403, 134),
    @(30, 5, "Murein Recycling", 0.9979720200255345, 0.1052631578947368, 38),
    @(31, 14, "Unassigned", 0.9980525965526649, 0.07692307692307693, 26),
    @(32, 15, "Transport, Inner Membrane", 0.9999933090549563, 0.1867469879518072, 332),
    @(33, 22, "Transport, Outer Membrane", 0.9999997977343006, 0.02173913043478261, 46),
    @(34, 31, "Alternate Carbon Metabolism", 0.9999999991307806, 0.1128205128205128, 195),
    @(35, 21, "Glycerophospholipid Metabolism", 0.9999999999999669, 0.0975609756097561, 246),
    @(36, 28, "Transport, Outer Membrane Porin", 0.999999999999999, 0.0962962962962963, 270),
    @(37, 33, "Methylglyoxal Metabolism", 1, 0, 9),
    @(38, 16, "Extracellular exchange", 1, 0.08333333333333333, 324),
    @(39, 34, "Nitrogen Metabolism", 1, 0, 13),
    @(40, 19, "tRNA Charging", 1, 0, 22)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(39, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(40, 1).PasteSpecial(-4122) | Out-Null

# ---- Sheet 3 data ----
$ws = $wb.Worksheets.Item(3)
$data = @(
    @(2, 17, "Cofactor and Prosthetic Group Biosynthesis", [double]"4.953145639386862e-36", 0.6575342465753424, 219),
    @(3, 25, "Membrane Lipid Metabolism", [double]"6.16701485270427e-15", 0.8536585365853658, 41),
    @(4, 6, "Tyrosine, Tryptophan, and Phenylalanine Metabolism", [double]"2.840218861388176e-09", 0.8695652173913043, 23),
    @(5, 0, "Purine and Pyrimidine Biosynthesis", [double]"4.021037266261487e-09", 0.84, 25),
    @(6, 30, "Valine, Leucine, and Isoleucine Metabolism", [double]"3.427793027790458e-08", 0.9375, 16),
    @(7, 3, "Pentose Phosphate Pathway", [double]"1.489779530763258e-07", 1, 12),
    @(8, 2, "Histidine Metabolism", [double]"2.067786117239401e-06", 1, 10),
    @(9, 38, "Arginine and Proline Metabolism", [double]"4.463464965068007e-05", 0.575, 40),
    @(10, 1, "Alanine and Aspartate Metabolism", 0.000195682403563015, 0.8888888888888888, 9),
    @(11, 27, "Murein Biosynthesis", 0.0002460644056356736, 0.7333333333333333, 15),
    @(12, 37, "Glycolysis/Gluconeogenesis", 0.0003476585209051788, 0.6363636363636364, 22),
    @(13, 35, "Citric Acid Cycle", 0.000678795666309634, 0.7142857142857143, 14),
    @(14, 10, "Lipopolysaccharide Biosynthesis / Recycling", 0.0007607742690504731, 0.4507042253521127, 71),
    @(15, 4, "Threonine and Lysine Metabolism", 0.001039269680393019, 0.631578947368421, 19),
    @(16, 36, "Methionine Metabolism", 0.001544258187657555, 0.6666666666666666, 15),
    @(17, 26, "Glyoxylate Metabolism", 0.005390887908528706, 1, 4),
    @(18, 23, "Pyruvate Metabolism", 0.005707166617169939, 0.7, 10),
    @(19, 7, "Cysteine Metabolism", 0.009486140483365409, 0.6153846153846154, 13),
    @(20, 9, "Glutamate Metabolism", 0.04983015248522543, 0.6666666666666666, 6),
    @(21, 13, "Intracellular demand", 0.04983015248522543, 0.6666666666666666, 6),
    @(22, 24, "Biomass and maintenance functions", 0.06353215463546867, 0.75, 4),
    @(23, 12, "Folate Metabolism", 0.06719145206329499, 0.5555555555555556, 9),
    @(24, 18, "Oxidative Phosphorylation", 0.1433694697850163, 0.3461538461538461, 52),
    @(25, 8, "Glycine and Serine Metabolism", 0.3233465839420369, 0.3571428571428572, 14),
    @(26, 11, "Anaplerotic Reactions", 0.3753073073139839, 0.375, 8),
    @(27, 33, "Methylglyoxal Metabolism", 0.4591845314337145, 0.3333333333333333, 9),
    @(28, 32, "Inorganic Ion Transport and Metabolism", 0.5713379358340012, 0.2678571428571428, 112),
    @(29, 20, "Cell Envelope Biosynthesis", 0.6406873664060805, 0.2611940298507462, 134),
    @(30, 14, "Unassigned", 0.9502516307631982, 0.1538461538461539, 26),
    @(31, 29, "Nucleotide Salvage Pathway", 0.9674232863632189, 0.2086330935251799, 139),
    @(32, 5, "Murein Recycling", 0.9993727594364125, 0.07894736842105263, 38),
    @(33, 22, "Transport, Outer Membrane", 0.9999995935156075, 0.02173913043478261, 46),
    @(34, 15, "Transport, Inner Membrane", 0.9999997886334522, 0.1626506024096386, 332),
    @(35, 31, "Alternate Carbon Metabolism", 0.9999998290577714, 0.1282051282051282, 195),
    @(36, 28, "Transport, Outer Membrane Porin", 0.9999999999989417, 0.1074074074074074, 270),
    @(37, 21, "Glycerophospholipid Metabolism", 0.999999999999526, 0.0975609756097561, 246),
    @(38, 16, "Extracellular exchange", 1, 0.08024691358024691, 324),
    @(39, 34, "Nitrogen Metabolism", 1, 0, 13),
    @(40, 19, "tRNA Charging", 1, 0, 22)
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
$ws.Cells.Item(2, 1).Copy() | Out-Null
$ws.Cells.Item(39, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(40, 1).PasteSpecial(-4122) | Out-Null

Write-Host "Edit complete"
